$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-18 Wednesday" "2025-06-19 Thursday"

Replace-Text "95×34=" "17×56="
Replace-Text "41×45=" "98×86="
Replace-Text "97×19=" "16×90="
Replace-Text "63×92=" "65×17="
Replace-Text "59×21=" "74×82="
Replace-Text "66×28=" "69×89="
Replace-Text "61×92=" "92×36="
Replace-Text "50×92=" "72×71="
Replace-Text "96×17=" "68×32="
Replace-Text "83×65=" "45×82="
Replace-Text "26×27=" "22×24="
Replace-Text "72×79=" "81×84="
Replace-Text "69×31=" "36×55="
Replace-Text "39×67=" "48×40="
Replace-Text "50×42=" "62×58="
Replace-Text "45×77=" "47×35="
Replace-Text "60×21=" "88×72="
Replace-Text "23×18=" "47×48="
Replace-Text "31×39=" "49×86="
Replace-Text "25×41=" "13×94="
Replace-Text "71×59=" "62×60="
Replace-Text "38×42=" "76×94="
Replace-Text "94×53=" "86×99="
Replace-Text "84×98=" "93×23="
Replace-Text "95×23=" "70×95="
